# Add a "Save" column (H) to the sheet, mirroring the style of the
# existing header cells (G1 "sum", etc.) and filling in the per-row
# save values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, styled like the other header cells (e.g. G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row "Save" values for H2:H17.
$saveValues = @(0, 1, 1, 0, 1, 0, 0, 0, 0, 1, 1, 0, 0, 0, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
